# Append three new days (2026-01-23, 2026-01-24, 2026-01-25) of daily
# charging-station data to Sheet1, two rows per day (四方坪站 / 高岭站),
# continuing directly after the existing last row (45).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A:日期(date serial)  B:站点  C:充电量(kwh)  D:充电总收入(元)  E:充电服务费收入(元)  F:总订单数量
$newRows = @(
    @(46045, "四方坪站", 14178.71,            11310.45,            4630.78,  614),
    @(46045, "高岭站",   5618.22,             4815.3599999999997, 1485.59,  209),
    @(46046, "四方坪站", 13637.41,            10897.33,            4421.6099999999997, 583),
    @(46046, "高岭站",   4735.83,             3785.78,             1271.5,   164),
    @(46047, "四方坪站", 10909.53,            8480.4500000000007, 3539.95,  462),
    @(46047, "高岭站",   5131.6899999999996,  4104.43,             1412.57,  171)
)

$startRow = 46

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rec = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rec[0]
    $ws.Cells.Item($r, 2).Value = $rec[1]
    $ws.Cells.Item($r, 3).Value = $rec[2]
    $ws.Cells.Item($r, 4).Value = $rec[3]
    $ws.Cells.Item($r, 5).Value = $rec[4]
    $ws.Cells.Item($r, 6).Value = $rec[5]
}

# Match the author's final scroll/selection state recorded in the sheet view.
$excel.ActiveWindow.ScrollRow = 38
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H48").Select()
